$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.8984962467368
$ws.Range("C2").Value = 11.24160509503961
$ws.Range("D2").Value = 5.969196063717503
$ws.Range("F2").Value = 57.31403531164443
$ws.Range("G2").Value = 3.7333297268994
$ws.Range("J2").Value = 10.7632301593864
$ws.Range("M2").Value = 19.46103087425163
$ws.Range("B3").Value = 18.64260506295506
$ws.Range("C3").Value = 11.0250145304494
$ws.Range("D3").Value = 5.936817003760138
$ws.Range("F3").Value = 56.35541699134136
$ws.Range("G3").Value = 3.738188842496593
$ws.Range("J3").Value = 10.77125333891725
$ws.Range("M3").Value = 19.46746815874981
$ws.Range("B4").Value = 18.49228409001401
$ws.Range("C4").Value = 10.89631912078327
$ws.Range("D4").Value = 5.917750549945874
$ws.Range("F4").Value = 55.76692645410509
$ws.Range("G4").Value = 3.741321186311438
$ws.Range("J4").Value = 10.77726643505485
$ws.Range("M4").Value = 19.47781167879008
$ws.Range("B5").Value = 18.43282339118874
$ws.Range("C5").Value = 10.84504907722565
$ws.Range("D5").Value = 5.910187050355404
$ws.Range("F5").Value = 55.52736877362982
$ws.Range("G5").Value = 3.742635235479936
$ws.Range("J5").Value = 10.77998980862089
$ws.Range("M5").Value = 19.48363130284492
$ws.Range("B6").Value = 18.42306093451035
$ws.Range("C6").Value = 10.83660938572837
$ws.Range("D6").Value = 5.908943628713014
$ws.Range("F6").Value = 55.48761267214643
$ws.Range("G6").Value = 3.742855707878715
$ws.Range("J6").Value = 10.78045850118826
$ws.Range("M6").Value = 19.48469446314044
$ws.Range("B7").Value = 18.49147479729257
$ws.Range("C7").Value = 10.8956227961106
$ws.Range("D7").Value = 5.917647709244327
$ws.Range("F7").Value = 55.76369434986085
$ws.Range("G7").Value = 3.741338755596524
$ws.Range("J7").Value = 10.77730205851547
$ws.Range("M7").Value = 19.47788367157538
$ws.Range("B8").Value = 18.80890780886543
$ws.Range("C8").Value = 11.16609156290284
$ws.Range("D8").Value = 5.957863691025628
$ws.Range("F8").Value = 56.98361993451871
$ws.Range("G8").Value = 3.734974363431582
$ws.Range("J8").Value = 10.76577080131484
$ws.Range("M8").Value = 19.46192281132033
$ws.Range("B9").Value = 19.48129604952349
$ws.Range("C9").Value = 11.72627421137064
$ws.Range("D9").Value = 6.043112924458401
$ws.Range("F9").Value = 59.36669242314252
$ws.Range("G9").Value = 3.723666925541245
$ws.Range("J9").Value = 10.75179543438049
$ws.Range("M9").Value = 19.48139776147003
$ws.Range("B10").Value = 20.00012227652342
$ws.Range("C10").Value = 12.15028923819359
$ws.Range("D10").Value = 6.109518921963804
$ws.Range("F10").Value = 61.09876512677802
$ws.Range("G10").Value = 3.716063627258094
$ws.Range("J10").Value = 10.7468127091906
$ws.Range("M10").Value = 19.52667720609165
$ws.Range("B11").Value = 20.24032207935586
$ws.Range("C11").Value = 12.34470708424749
$ws.Range("D11").Value = 6.140516597235359
$ws.Range("F11").Value = 61.87997365874671
$ws.Range("G11").Value = 3.712755275750523
$ws.Range("J11").Value = 10.7456972990669
$ws.Range("M11").Value = 19.55398049136167
$ws.Range("B12").Value = 20.3317796068098
$ws.Range("C12").Value = 12.41845419541169
$ws.Range("D12").Value = 6.152364938653754
$ws.Range("F12").Value = 62.17461620349076
$ws.Range("G12").Value = 3.711523938004527
$ws.Range("J12").Value = 10.74544070820899
$ws.Range("M12").Value = 19.56527998973617
$ws.Range("B13").Value = 20.31206201059521
$ws.Range("C13").Value = 12.40256732916072
$ws.Range("D13").Value = 6.149808340556933
$ws.Range("F13").Value = 62.11121569827633
$ws.Range("G13").Value = 3.711788176716142
$ws.Range("J13").Value = 10.74548859249611
$ws.Range("M13").Value = 19.56280380422424
$ws.Range("B14").Value = 20.24783692455953
$ws.Range("C14").Value = 12.35077232104114
$ws.Range("D14").Value = 6.141489178722001
$ws.Range("F14").Value = 61.90423902650234
$ws.Range("G14").Value = 3.712653543662945
$ws.Range("J14").Value = 10.74567286519998
$ws.Range("M14").Value = 19.55489088794382
$ws.Range("B15").Value = 20.20855918316187
$ws.Range("C15").Value = 12.31905995461051
$ws.Range("D15").Value = 6.136407700338501
$ws.Range("F15").Value = 61.77729931270135
$ws.Range("G15").Value = 3.71318639605099
$ws.Range("J15").Value = 10.74580733511483
$ws.Range("M15").Value = 19.55016892610217
$ws.Range("B16").Value = 19.98450004150408
$ws.Range("C16").Value = 12.13760618187747
$ws.Range("D16").Value = 6.107508756488635
$ws.Range("F16").Value = 61.04755726404403
$ws.Range("G16").Value = 3.716282848424475
$ws.Range("J16").Value = 10.74690879134059
$ws.Range("M16").Value = 19.52502749729423
$ws.Range("B17").Value = 19.84804513188681
$ws.Range("C17").Value = 12.02661356905648
$ws.Range("D17").Value = 6.089980070346619
$ws.Range("F17").Value = 60.59800522262335
$ws.Range("G17").Value = 3.718220830109018
$ws.Range("J17").Value = 10.7478795437306
$ws.Range("M17").Value = 19.51131908113452
$ws.Range("B18").Value = 19.76995969276845
$ws.Range("C18").Value = 11.96292312897134
$ws.Range("D18").Value = 6.079972536177348
$ws.Range("F18").Value = 60.33881475940195
$ws.Range("G18").Value = 3.719349674881136
$ws.Range("J18").Value = 10.74854624383027
$ws.Range("M18").Value = 19.50406582124171
$ws.Range("B19").Value = 19.74359311620425
$ws.Range("C19").Value = 11.94138728813667
$ws.Range("D19").Value = 6.076597049340301
$ws.Range("F19").Value = 60.25095761484599
$ws.Range("G19").Value = 3.719734320972624
$ws.Range("J19").Value = 10.74879057827489
$ws.Range("M19").Value = 19.50171854443964
$ws.Range("B20").Value = 19.86253035710077
$ws.Range("C20").Value = 12.03841414231548
$ws.Range("D20").Value = 6.091838341292502
$ws.Range("F20").Value = 60.64592658076574
$ws.Range("G20").Value = 3.718013063433503
$ws.Range("J20").Value = 10.74776498994005
$ws.Range("M20").Value = 19.51271303627983
$ws.Range("B21").Value = 20.2666886535658
$ws.Range("C21").Value = 12.36598310550451
$ws.Range("D21").Value = 6.143929755080173
$ws.Range("F21").Value = 61.96506692161488
$ws.Range("G21").Value = 3.712398783249629
$ws.Range("J21").Value = 10.74561423853544
$ws.Range("M21").Value = 19.55718907704798
$ws.Range("B22").Value = 20.53369215423791
$ws.Range("C22").Value = 12.58075774201812
$ws.Range("D22").Value = 6.17861510424164
$ws.Range("F22").Value = 62.82020462427217
$ws.Range("G22").Value = 3.708854552226083
$ws.Range("J22").Value = 10.74517504676261
$ws.Range("M22").Value = 19.59185166042861
$ws.Range("B23").Value = 20.39095885034827
$ws.Range("C23").Value = 12.46609563371426
$ws.Range("D23").Value = 6.160045418216804
$ws.Range("F23").Value = 62.3645113412674
$ws.Range("G23").Value = 3.710734790749841
$ws.Range("J23").Value = 10.74532094866134
$ws.Range("M23").Value = 19.57284125083687
$ws.Range("B24").Value = 19.85598044426064
$ws.Range("C24").Value = 12.03307871996969
$ws.Range("D24").Value = 6.090997998899941
$ws.Range("F24").Value = 60.62426361373545
$ws.Range("G24").Value = 3.718106949030872
$ws.Range("J24").Value = 10.7478164414575
$ws.Range("M24").Value = 19.51208087313676
$ws.Range("B25").Value = 19.29462575419445
$ws.Range("C25").Value = 11.57213438199699
$ws.Range("D25").Value = 6.019378230269244
$ws.Range("F25").Value = 58.72437094659625
$ws.Range("G25").Value = 3.726601422083518
$ws.Range("J25").Value = 10.75464934603525
$ws.Range("M25").Value = 19.47068701852979
